$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.098603367805481
$ws.Range("B1").Value = 3.831497430801392
$ws.Range("C1").Value = 3.616532325744629
$ws.Range("D1").Value = 3.262124061584473
$ws.Range("E1").Value = 1.235480070114136
